# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (374, 375) into the Pomelo sheet,
# pushing the existing data rows (previously 374:472) down to 376:474.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 374; this shifts the
# existing rows 374:472 down to 376:474, preserving their contents/format.
$ws.Rows("374:375").Insert()

# --- New row 374 ---
$ws.Range("A374").Value = 10
$ws.Range("B374").Value = "Vega Modelo de Temuco"
$ws.Range("C374").Value = "La Araucanía"
$ws.Range("D374").Value = 45135
$ws.Range("E374").Value = 9
$ws.Range("F374").Value = "Fruta"
$ws.Range("G374").Value = 100102
$ws.Range("H374").Value = "Cítricos"
$ws.Range("I374").Value = 100102006
$ws.Range("J374").Value = "Pomelo"
$ws.Range("K374").Value = "Start Ruby"
$ws.Range("L374").Value = "Primera"
$ws.Range("M374").Value = 100
$ws.Range("N374").Value = 15000
$ws.Range("O374").Value = 15000
$ws.Range("P374").Value = 15000
$ws.Range("Q374").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R374").Value = "Región de O'Higgins"
$ws.Range("S374").Value = 1000
$ws.Range("T374").Value = 15

# --- New row 375 ---
$ws.Range("A375").Value = 10
$ws.Range("B375").Value = "Vega Modelo de Temuco"
$ws.Range("C375").Value = "La Araucanía"
$ws.Range("D375").Value = 45135
$ws.Range("E375").Value = 9
$ws.Range("F375").Value = "Fruta"
$ws.Range("G375").Value = 100102
$ws.Range("H375").Value = "Cítricos"
$ws.Range("I375").Value = 100102006
$ws.Range("J375").Value = "Pomelo"
$ws.Range("K375").Value = "Start Ruby"
$ws.Range("L375").Value = "Primera"
$ws.Range("M375").Value = 1
$ws.Range("N375").Value = 280000
$ws.Range("O375").Value = 280000
$ws.Range("P375").Value = 280000
$ws.Range("Q375").Value = "`$/bins (350 kilos)"
$ws.Range("R375").Value = "Región de O'Higgins"
$ws.Range("S375").Value = 800
$ws.Range("T375").Value = 350
